$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it to make the edits, then restore protection.
$ws.Unprotect("D382")

# Update the confidential disclaimer text (shared string used by A16): the "as of" date
# changes from 2021-05-14 to 2021-05-17.
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."
$ws.Rows.Item(16).AutoFit()

# Update the Weight (column D) and Percent Change (column E) values for rows 2-13.
$ws.Range("D2").Value = 0.03068785630194209
$ws.Range("E2").Value = 0.0006882312456986117
$ws.Range("D3").Value = 0.02377355481689732
$ws.Range("E3").Value = -0.008913910391742941
$ws.Range("D4").Value = 0.05235614021723576
$ws.Range("E4").Value = -0.001618122977346093
$ws.Range("D5").Value = 0.1375133574295073
$ws.Range("E5").Value = -0.002120026092628846
$ws.Range("D6").Value = 0.03212854836570085
$ws.Range("E6").Value = 0.024171270718232
$ws.Range("D7").Value = 0.1173423633472372
$ws.Range("E7").Value = -0.0008989341209707513
$ws.Range("D8").Value = 0.1034286857428439
$ws.Range("E8").Value = -0.003081384810585508
$ws.Range("D9").Value = 0.02978466506579223
$ws.Range("E9").Value = 0.009752133279155029
$ws.Range("D10").Value = 0.1281221619541014
$ws.Range("E10").Value = 0.00166697536580851
$ws.Range("D11").Value = 0.2424627077654275
$ws.Range("E11").Value = -0.006880520688052005
$ws.Range("D12").Value = 0.1023999589933143
$ws.Range("E12").Value = -0.008759976640062339
$ws.Range("E13").Value = -0.002275894784628307

# Restore sheet protection (password + locked contents) as it was before the edit.
$ws.Protect("D382", $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
